$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "692×6=" "692×2="
Replace-Text "413×8=" "305×2="
Replace-Text "369×5=" "128×6="
Replace-Text "765×6=" "762×6="
Replace-Text "427×6=" "374×8="
Replace-Text "153×2=" "856×2="
Replace-Text "297×9=" "502×2="
Replace-Text "728×4=" "753×3="
Replace-Text "712×7=" "237×7="
Replace-Text "885×9=" "414×4="
Replace-Text "446×3=" "722×6="
Replace-Text "652×6=" "705×8="
Replace-Text "803×7=" "473×6="
Replace-Text "151×9=" "702×6="
Replace-Text "731×8=" "308×5="
Replace-Text "209×9=" "754×8="
Replace-Text "438×8=" "888×2="
Replace-Text "571×9=" "739×5="
Replace-Text "320×2=" "772×2="
Replace-Text "898×5=" "863×3="
Replace-Text "578×2=" "388×2="
Replace-Text "491×7=" "613×8="
Replace-Text "234×2=" "616×8="
Replace-Text "646×9=" "342×6="
Replace-Text "818×5=" "456×3="
